$d = $word.ActiveDocument

# Find the paragraph containing the "missing feature access or service call"
# error message that follows "Expression ""self."" is invalid:" text, and
# position the insertion point right after that run (i.e. at the end of the
# paragraph, just before the paragraph mark).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like '*Expression*self*is invalid*missing feature access or service call*') {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

$insertPos = $target.Range.End - 1  # just before the paragraph mark

function Insert-ErrorRun([int]$pos, [string]$text, [bool]$styled) {
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text)
    $endPos = $pos + $text.Length
    if ($styled) {
        $fr = $d.Range($pos, $endPos)
        $fr.Font.Color = 255
        $fr.Font.Size = 16
        $fr.Font.HighlightColorIndex = 16
    }
    return $endPos
}

$pos = $insertPos
$pos = Insert-ErrorRun $pos "    " $false
$pos = Insert-ErrorRun $pos "<---" $true
$pos = Insert-ErrorRun $pos "Couldn't find the 'self' variable" $true
$pos = Insert-ErrorRun $pos "    " $false
$pos = Insert-ErrorRun $pos "<---" $true
$pos = Insert-ErrorRun $pos "missing feature access or service call" $true
